$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

foreach ($r in 4..7) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-13 02:36:14"

    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-13 02:36:22"

    $overview.Range("G$r").Value = "2016-08-13 02:36:22"
}
